# Insert a new record row at row 87 (shifts the existing rows 87-141 down to 88-142)
# and populate it with the new price-observation data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("87:87").Insert()

$ws.Cells.Item(87, 1).Value = 3
$ws.Cells.Item(87, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(87, 3).Value = "Coquimbo"
$ws.Cells.Item(87, 4).Value = 44606
$ws.Cells.Item(87, 5).Value = 5
$ws.Cells.Item(87, 6).Value = 100112030
$ws.Cells.Item(87, 7).Value = "Poroto granado"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 73
$ws.Cells.Item(87, 11).Value = 24000
$ws.Cells.Item(87, 12).Value = 25000
$ws.Cells.Item(87, 13).Value = 24479
$ws.Cells.Item(87, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(87, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(87, 16).Value = 979
$ws.Cells.Item(87, 17).Value = 25
$ws.Cells.Item(87, 18).Value = "Hortaliza"
